# Apply the "descriptor builder" edits to the Batting worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Batting")

# --- 1. Rename the "fancy" aspect to "utility_hitter" (row 1, column E) ---
$ws.Range("E1").Value = "utility_hitter"

# --- 2. Update the smallball legend note (Y2) ---
$ws.Range("Y2").Value = "smallball: contact, control, speed, discipline"

# --- 3. Row 6 ("slugging" sub-row under the "slugging" aspect): rename labels ---
$ws.Range("C6").Value = "Slugger Only"
$ws.Range("D6").Value = "Pure Slugger"

# --- 4. Row 12 ("slugging" sub-row under the "smallball" aspect):
#         replace computed formulas with literal overrides ---
$ws.Range("C12").Value = "Contact Only"
$ws.Range("D12").Value = "Reliable Contact"
$ws.Range("E12").Value = "Omnipresent Contact"

# --- 5. Rows 24-27 (the "utility_hitter"/old "fancy" aspect group):
#         rename the "Weak Strategy Hitter" label to "Weak Utility" ---
$ws.Range("C24").Value = "Weak Utility"
$ws.Range("C25").Value = "Weak Utility"
$ws.Range("C26").Value = "Weak Utility"
$ws.Range("C27").Value = "Weak Utility"

# --- 6. Fix the M-column generator formulas so each emitted python list literal
#         ends with a trailing comma (']' -> '],') ---
$m5 = "=IF(B5=`"`",`nIF(B6=`"`",IF(B4=`"`",`"`",`"    },`"),`nCONCATENATE(`"    `",`"'`",A5,`"': {`")),`nCONCATENATE(`"        `",`"'`",B5,`"': [`",IF(C5=`"`",`"`",`"'`"&C5&`"',`"),IF(D5=`"`",`"`",`" '`"&D5&`"',`"),IF(E5=`"`",`"`",`" '`"&E5&`"',`"),IF(F5=`"`",`"`",`" '`"&F5&`"',`"),IF(G5=`"`",`"`",`" '`"&G5&`"',`"),IF(H5=`"`",`"`",`" '`"&H5&`"',`"),IF(I5=`"`",`"`",`" '`"&I5&`"',`"),IF(J5=`"`",`"`",`" '`"&J5&`"',`"),IF(K5=`"`",`"`",`" '`"&K5&`"',`"),IF(L5=`"`",`"`",`" '`"&L5&`"',`"),`"],`"))"
$ws.Range("M5").Formula = $m5

$m6 = "=IF(B6=`"`",`nIF(B7=`"`",IF(B5=`"`",`"`",`"    },`"),`nCONCATENATE(`"    `",`"'`",A6,`"': {`")),`nCONCATENATE(`"        `",`"'`",B6,`"': [`",IF(C6=`"`",`"`",`"'`"&C6&`"',`"),IF(D6=`"`",`"`",`" '`"&D6&`"',`"),IF(E6=`"`",`"`",`" '`"&E6&`"',`"),IF(F6=`"`",`"`",`" '`"&F6&`"',`"),IF(G6=`"`",`"`",`" '`"&G6&`"',`"),IF(H6=`"`",`"`",`" '`"&H6&`"',`"),IF(I6=`"`",`"`",`" '`"&I6&`"',`"),IF(J6=`"`",`"`",`" '`"&J6&`"',`"),IF(K6=`"`",`"`",`" '`"&K6&`"',`"),IF(L6=`"`",`"`",`" '`"&L6&`"',`"),`"],`"))"
$ws.Range("M6:M69").Formula = $m6

$m70 = "=IF(B70=`"`",`nIF(B71=`"`",IF(B69=`"`",`"`",`"    },`"),`nCONCATENATE(`"    `",`"'`",A70,`"': {`")),`nCONCATENATE(`"        `",`"'`",B70,`"': [`",IF(C70=`"`",`"`",`"'`"&C70&`"',`"),IF(D70=`"`",`"`",`" '`"&D70&`"',`"),IF(E70=`"`",`"`",`" '`"&E70&`"',`"),IF(F70=`"`",`"`",`" '`"&F70&`"',`"),IF(G70=`"`",`"`",`" '`"&G70&`"',`"),IF(H70=`"`",`"`",`" '`"&H70&`"',`"),IF(I70=`"`",`"`",`" '`"&I70&`"',`"),IF(J70=`"`",`"`",`" '`"&J70&`"',`"),IF(K70=`"`",`"`",`" '`"&K70&`"',`"),IF(L70=`"`",`"`",`" '`"&L70&`"',`"),`"],`"))"
$ws.Range("M70:M120").Formula = $m70

# --- 7. Restore the active-cell selection recorded in the workbook ---
$ws.Range("D13").Select()

Write-Host "Edits applied"
